# Improve cell segmentation algorithm: use histogram for identifying cells
# with no inclusions. Updates the Number_of_Inclusions (col B) and the
# derived Number_of_Inclusions_per_Nucleus (col D = B / C) for the affected
# image rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 9
    3  = 2
    4  = 3
    5  = 3
    6  = 6
    7  = 1
    8  = 2
    9  = 10
    10 = 30
    11 = 4
    13 = 5
    14 = 5
    15 = 22
    17 = 0
    19 = 9
    21 = 5
    22 = 8
    23 = 4
    24 = 16
    25 = 8
    26 = 41
    27 = 2
    28 = 0
    29 = 16
    30 = 1
    31 = 6
    32 = 10
    33 = 2
    34 = 1
    36 = 1
    37 = 3
    38 = 2
    40 = 6
    43 = 1
    45 = 8
    46 = 1
    47 = 4
    48 = 8
    51 = 3
    52 = 12
    53 = 1
    54 = 0
    58 = 9
    60 = 11
}

foreach ($row in $updates.Keys) {
    $newB = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $newB
    $nucCount = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $newB / $nucCount
}
